# Commit: "Fruta / hortaliza, semanal"
# A new weekly price-record row is inserted at row 464 (shifting the existing
# rows 464-521 down to 465-522), and the new row 464 is populated with the
# latest observation for Berenjena / Vega Modelo de Temuco.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 464; this pushes the old 464..521
# rows down to 465..522 and copies formatting (incl. the date style) from
# the row above, matching the target dimension A1:R522.
$ws.Rows("464:464").Insert()

# Populate the newly inserted row 464 with the new weekly data point.
$ws.Cells.Item(464, 1).Value = 10
$ws.Cells.Item(464, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(464, 3).Value = "La Araucanía"
$ws.Cells.Item(464, 4).Value = 45212
$ws.Cells.Item(464, 5).Value = 9
$ws.Cells.Item(464, 6).Value = 100112001
$ws.Cells.Item(464, 7).Value = "Berenjena"
$ws.Cells.Item(464, 8).Value = "Sin especificar"
$ws.Cells.Item(464, 9).Value = "Primera"
$ws.Cells.Item(464, 10).Value = 40
$ws.Cells.Item(464, 11).Value = 12000
$ws.Cells.Item(464, 12).Value = 12000
$ws.Cells.Item(464, 13).Value = 12000
$ws.Cells.Item(464, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(464, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(464, 16).Value = 300
$ws.Cells.Item(464, 17).Value = 40
$ws.Cells.Item(464, 18).Value = "Hortaliza"
